$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the idiom text in B18 ("Пробел в знаниях" -> "Двойное дно")
$ws.Range("B18").Value = "Двойное дно"

# Update the selected cell to match the saved view state
$ws.Range("B19").Select()
